# Insert a new weekly price row before row 6 (shifting the old row 6 data
# down to row 7), then update row 6 with the latest figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 6 down to row 7 by inserting a new blank row at position 6.
$ws.Rows.Item(6).Insert()

# Row 6: new weekly entry (values taken from the diff).
# (NumberFormat for column D is inherited from the row above by Insert(),
# so no explicit style/number-format assignment is needed here.)
$ws.Cells.Item(6, 1).Value = 10
$ws.Cells.Item(6, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(6, 3).Value = "La Araucanía"
$ws.Cells.Item(6, 4).Value = 44474
$ws.Cells.Item(6, 5).Value = 9
$ws.Cells.Item(6, 6).Value = 100112042
$ws.Cells.Item(6, 7).Value = "Locoto"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 20
$ws.Cells.Item(6, 11).Value = 1600
$ws.Cells.Item(6, 12).Value = 1600
$ws.Cells.Item(6, 13).Value = 1600
$ws.Cells.Item(6, 14).Value = "$/kilo"
$ws.Cells.Item(6, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(6, 16).Value = 1600
$ws.Cells.Item(6, 17).Value = 1
$ws.Cells.Item(6, 18).Value = "Hortaliza"
